$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 with the merged tuple-style text, then delete the now-unused rows 3 and 4
$ws.Range("A2").Value = "('Elemental Shaman', ['Token Creature — Elemental Shaman', '3/1'])"

$ws.Range("A3").EntireRow.Delete()
$ws.Range("A3").EntireRow.Delete()
